# Planification.xlsx - apply the "Cahier des charges - conclusion" update:
#  - AGENDA: move the "Rencontre prof. 17h / DEBUT SPRINT 0" note from B6 to
#    B4 (the actual Sprint-0 kickoff date), and leave B6 blank again.
#  - SPRINTS: the Sprint 0 start date moves from 2020-02-24 to 2020-02-18.
#  - New worksheet "SPRINT 0" tracking the hours already logged against the
#    cahier des charges work.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# AGENDA sheet: shift the "Rencontre prof." annotation up from row 6 to
# row 4, since the professor meeting / sprint-0 kickoff actually happened
# on 2020-02-18 (row 4 date range), not 2020-02-24 (row 6 date range).
# ---------------------------------------------------------------------
$agenda = $wb.Worksheets.Item("AGENDA")

# Give B4 the same "note" formatting B6 currently has, then write the text.
$agenda.Range("B6").Copy()
$agenda.Range("B4").PasteSpecial(-4122) # xlPasteFormats
$agenda.Range("B4").Value = "Rencontre prof. 17h" + [char]10 + "DEBUT SPRINT 0"

# Restore B6 back to the plain unfilled note style and clear its text.
$agenda.Range("A4").Copy()
$agenda.Range("B6").PasteSpecial(-4122) # xlPasteFormats
$agenda.Range("B6").ClearContents()

$agenda.Range("E15").Select() | Out-Null

# ---------------------------------------------------------------------
# SPRINTS sheet: Sprint 0 actually started 2020-02-18, not 2020-02-24.
# ---------------------------------------------------------------------
$sprints = $wb.Worksheets.Item("SPRINTS")
$sprints.Range("B2").Value = 43879
$sprints.Range("B3").Select() | Out-Null

# ---------------------------------------------------------------------
# New "SPRINT 0" worksheet: hours logged on the cahier des charges.
# ---------------------------------------------------------------------
$sprint0 = $wb.Worksheets.Add($null, $sprints)
$sprint0.Name = "SPRINT 0"

$day1 = Get-Date -Year 2020 -Month 2 -Day 18 -Hour 0 -Minute 0 -Second 0
$day2 = Get-Date -Year 2020 -Month 2 -Day 24 -Hour 0 -Minute 0 -Second 0

# Log the entries first (matches the order the rows were actually typed),
# then add the column headers afterwards.
$sprint0.Range("B3").Value = $day1
$sprint0.Range("C3").Value = "Rencontre professeur, discussions des objectifs du travail"
$sprint0.Range("D3").Value = 1

$sprint0.Range("B2").Value = "Date"
$sprint0.Range("C2").Value = "Quoi"
$sprint0.Range("D2").Value = "Temps (h)"

$sprint0.Range("B4").Value = $day2
$sprint0.Range("C4").Value = "Introduction du cahier des charges"
$sprint0.Range("D4").Value = 1

$sprint0.Range("B5").Value = $day2
$sprint0.Range("C5").Value = "Mise en page cahier des charges"
$sprint0.Range("D5").Value = 0.5

$sprint0.Range("B6").Value = $day2
$sprint0.Range("C6").Value = "Conclusion du cahier des charges"
$sprint0.Range("D6").Value = 0.5

$sprint0.Columns.Item(3).AutoFit()

$sprint0.Range("E6").Select() | Out-Null
